$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings) ---
# A8: "Volume 29   Number  42" -> "...43"  (last run's text, chars 21-22)
$ws.Range("A8").Characters(21, 2).Text = "43"

# C9: "Report Covering the Week  10/17/2022  Through  10/23/2022"
#     -> "...10/24/2022  Through  10/30/2022"
$ws.Range("C9").Characters(27, 10).Text = "10/24/2022"
$ws.Range("C9").Characters(48, 10).Text = "10/30/2022"

# --- Row 14 (Murder): C14 becomes a text "0" cell (like D14/E14) ---
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C14").PasteSpecial(-4122)

# --- Row 30 (Hate Crimes): C30 becomes text "0"; D30/E30 become numbers ---
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C30").PasteSpecial(-4122)

$ws.Range("F30").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value2 = 1

$ws.Range("H30").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value2 = -100

# --- Remaining numeric cell updates across rows 14-30 ---
$ws.Range("F14").Value2 = 1
$ws.Range("H14").Value2 = 0
$ws.Range("M14").Value2 = -8.333333333333
$ws.Range("C15").Value2 = 10
$ws.Range("D15").Value2 = 4
$ws.Range("E15").Value2 = 150
$ws.Range("F15").Value2 = 20
$ws.Range("G15").Value2 = 14
$ws.Range("H15").Value2 = 42.857142857142
$ws.Range("I15").Value2 = 174
$ws.Range("J15").Value2 = 137
$ws.Range("K15").Value2 = 27.007299270073
$ws.Range("L15").Value2 = 41.463414634146
$ws.Range("M15").Value2 = 52.631578947368
$ws.Range("N15").Value2 = 4.191616766467
$ws.Range("C16").Value2 = 61
$ws.Range("D16").Value2 = 35
$ws.Range("E16").Value2 = 74.285714285714
$ws.Range("F16").Value2 = 170
$ws.Range("G16").Value2 = 104
$ws.Range("H16").Value2 = 63.461538461538
$ws.Range("I16").Value2 = 1532
$ws.Range("J16").Value2 = 965
$ws.Range("K16").Value2 = 58.756476683937
$ws.Range("L16").Value2 = 48.737864077669
$ws.Range("M16").Value2 = -3.404791929382
$ws.Range("N16").Value2 = -78.973373593192
$ws.Range("C17").Value2 = 50
$ws.Range("D17").Value2 = 31
$ws.Range("E17").Value2 = 61.290322580645
$ws.Range("F17").Value2 = 198
$ws.Range("G17").Value2 = 176
$ws.Range("H17").Value2 = 12.5
$ws.Range("I17").Value2 = 2065
$ws.Range("J17").Value2 = 1666
$ws.Range("K17").Value2 = 23.949579831932
$ws.Range("L17").Value2 = 29.874213836478
$ws.Range("M17").Value2 = 61.833855799373
$ws.Range("N17").Value2 = -15.988608624898
$ws.Range("C18").Value2 = 38
$ws.Range("D18").Value2 = 36
$ws.Range("E18").Value2 = 5.555555555555
$ws.Range("F18").Value2 = 159
$ws.Range("G18").Value2 = 129
$ws.Range("H18").Value2 = 23.255813953488
$ws.Range("I18").Value2 = 1628
$ws.Range("J18").Value2 = 1320
$ws.Range("K18").Value2 = 23.333333333333
$ws.Range("L18").Value2 = 8.029197080291
$ws.Range("M18").Value2 = -25.763793889648
$ws.Range("N18").Value2 = -86.848695371193
$ws.Range("C19").Value2 = 145
$ws.Range("D19").Value2 = 118
$ws.Range("E19").Value2 = 22.881355932203
$ws.Range("F19").Value2 = 581
$ws.Range("G19").Value2 = 434
$ws.Range("H19").Value2 = 33.870967741935
$ws.Range("I19").Value2 = 5882
$ws.Range("J19").Value2 = 3536
$ws.Range("K19").Value2 = 66.346153846153
$ws.Range("L19").Value2 = 65.039281705948
$ws.Range("M19").Value2 = 75.897129186602
$ws.Range("N19").Value2 = -14.840017373678
$ws.Range("C20").Value2 = 43
$ws.Range("D20").Value2 = 26
$ws.Range("E20").Value2 = 65.384615384615
$ws.Range("F20").Value2 = 196
$ws.Range("G20").Value2 = 117
$ws.Range("H20").Value2 = 67.521367521367
$ws.Range("I20").Value2 = 1591
$ws.Range("J20").Value2 = 1143
$ws.Range("K20").Value2 = 39.195100612423
$ws.Range("L20").Value2 = 56.440511307767
$ws.Range("M20").Value2 = 8.157715839564
$ws.Range("N20").Value2 = -91.936955199675
$ws.Range("C21").Value2 = 347
$ws.Range("D21").Value2 = 250
$ws.Range("E21").Value2 = 38.8
$ws.Range("F21").Value2 = 1325
$ws.Range("G21").Value2 = 975
$ws.Range("H21").Value2 = 35.897435897435
$ws.Range("I21").Value2 = 12894
$ws.Range("J21").Value2 = 8792
$ws.Range("K21").Value2 = 46.656050955414
$ws.Range("L21").Value2 = 45.579767415603
$ws.Range("M21").Value2 = 28.836930455635
$ws.Range("N21").Value2 = -73.703424224501
$ws.Range("C22").Value2 = 6
$ws.Range("E22").Value2 = 100
$ws.Range("F22").Value2 = 31
$ws.Range("H22").Value2 = 158.333333333333
$ws.Range("I22").Value2 = 254
$ws.Range("J22").Value2 = 127
$ws.Range("L22").Value2 = 149.019607843137
$ws.Range("M22").Value2 = 58.75
$ws.Range("D23").Value2 = 3
$ws.Range("E23").Value2 = 66.666666666666
$ws.Range("F23").Value2 = 19
$ws.Range("G23").Value2 = 22
$ws.Range("H23").Value2 = -13.636363636363
$ws.Range("I23").Value2 = 196
$ws.Range("J23").Value2 = 188
$ws.Range("K23").Value2 = 4.255319148936
$ws.Range("L23").Value2 = 7.692307692307
$ws.Range("M23").Value2 = 38.028169014084
$ws.Range("C24").Value2 = 292
$ws.Range("D24").Value2 = 224
$ws.Range("E24").Value2 = 30.357142857142
$ws.Range("F24").Value2 = 1159
$ws.Range("G24").Value2 = 1092
$ws.Range("H24").Value2 = 6.135531135531
$ws.Range("I24").Value2 = 12692
$ws.Range("J24").Value2 = 9952
$ws.Range("K24").Value2 = 27.532154340836
$ws.Range("L24").Value2 = 50.539675008895
$ws.Range("M24").Value2 = 60.637893937476
$ws.Range("C25").Value2 = 100
$ws.Range("D25").Value2 = 79
$ws.Range("E25").Value2 = 26.582278481012
$ws.Range("F25").Value2 = 420
$ws.Range("G25").Value2 = 386
$ws.Range("H25").Value2 = 8.808290155440
$ws.Range("I25").Value2 = 4146
$ws.Range("J25").Value2 = 3631
$ws.Range("K25").Value2 = 14.183420545304
$ws.Range("L25").Value2 = 31.202531645569
$ws.Range("M25").Value2 = 3.443113772455
$ws.Range("C26").Value2 = 12
$ws.Range("D26").Value2 = 7
$ws.Range("E26").Value2 = 71.428571428571
$ws.Range("F26").Value2 = 24
$ws.Range("G26").Value2 = 24
$ws.Range("H26").Value2 = 0
$ws.Range("I26").Value2 = 248
$ws.Range("J26").Value2 = 226
$ws.Range("K26").Value2 = 9.734513274336
$ws.Range("L26").Value2 = 18.095238095238
$ws.Range("C27").Value2 = 11
$ws.Range("D27").Value2 = 18
$ws.Range("E27").Value2 = -38.888888888888
$ws.Range("G27").Value2 = 46
$ws.Range("H27").Value2 = 15.217391304347
$ws.Range("I27").Value2 = 529
$ws.Range("J27").Value2 = 464
$ws.Range("K27").Value2 = 14.008620689655
$ws.Range("L27").Value2 = 23.887587822014
$ws.Range("G28").Value2 = 4
$ws.Range("H28").Value2 = -50
$ws.Range("L28").Value2 = -14.864864864864
$ws.Range("M28").Value2 = 61.538461538461
$ws.Range("N28").Value2 = -71.493212669683
$ws.Range("G29").Value2 = 3
$ws.Range("H29").Value2 = -33.333333333333
$ws.Range("L29").Value2 = -8.620689655172
$ws.Range("M29").Value2 = 60.606060606060
$ws.Range("N29").Value2 = -72.959183673469
$ws.Range("G30").Value2 = 2
$ws.Range("H30").Value2 = 100
$ws.Range("J30").Value2 = 65
$ws.Range("K30").Value2 = -16.923076923076
